$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# The long product-name string (shared string) is shown in B1 on both sheets.
# Update its text on both sheets so the underlying shared string is edited in place.
$newProductName = "4238-RBI-EI-DB-SAR-REC-RNI-FEE+INT-FFConMONTonDAY25-FIFC-1-FFROP-DL-FIFR-1-MD-TR-1-OT-PER-1st"
$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Change the "shortname" cell B2 on ProductLoanInput from the numeric 4238
# to the text value "423p".
$wsInput.Range("B2").Value = "423p"

# Update the remembered cell selection on ProductLoanInput to B7.
$wsInput.Range("B7").Select()

# Make ProductLoanOutput the active sheet/tab, with B1 selected there.
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
